$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the Detailed Summary Table: D9 (Mobilization / PSUV relocation votes)
# changes from 92210 to 92105. All dependent formulas (E9, F9, F10, E11, D12,
# F12, G12, J25, K25, L25, J26, K26) recalculate automatically.
$ws.Range("D9").Value = 92105

# Update the active selection to F24, matching the saved view state.
$ws.Range("F24").Select()

$excel.Calculate()
